# Add a new FAQ entry ("time-page") with 5 question variants about the
# page's opening hours. The new rows are inserted at row 12, pushing the
# existing rows down by 5 (old rows 12-100 become rows 17-105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current row 12 ("invite-eating" block).
$ws.Rows("12:16").Insert()

# Fill in the new "time-page" tag rows.
$ws.Range("A12").Value = "time-page"
$ws.Range("B12").Value = "เวลาทำการของเพจ"

$ws.Range("A13").Value = "time-page"
$ws.Range("B13").Value = "เวลาของเพจ"

$ws.Range("A14").Value = "time-page"
$ws.Range("B14").Value = "เพจเปิดตอนไหน"

$ws.Range("A15").Value = "time-page"
$ws.Range("B15").Value = "เพจเปิดตีเท่าไหร่"

$ws.Range("A16").Value = "time-page"
$ws.Range("B16").Value = "เวลาเพจเปิด"

# Give column A an explicit width (matches the stored OOXML width="12").
$ws.Columns("A").ColumnWidth = 11.285714285714286

# Update the view: scroll back to the top and move the selection to B17
# (the row that used to be row 12 before the insert).
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("B17").Select()
